$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: Rule 1 text - "25x25Kms" -> "15x15Kms"
# Rebuilt as several runs (mirroring the way the original author typed
# the two digit corrections one at a time) by injecting OOXML directly
# into the run-level range of the "Rule 1" paragraph. This keeps the
# paragraph's own formatting (<w:pPr> indent) untouched because the
# paragraph mark itself is excluded from the replaced range.
# ---------------------------------------------------------------------
$rule1ParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*25x25Kms*") {
        $rule1ParaIndex = $i
    }
}

$pRule1 = $d.Paragraphs($rule1ParaIndex)
$runRange = $d.Range($pRule1.Range.Start, $pRule1.Range.End - 1)

$rule1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r><w:t xml:space="preserve">Rule 1: The city size should be limited to a maximum of </w:t></w:r>' + `
  '<w:r><w:t>1</w:t></w:r>' + `
  '<w:r><w:t>5x</w:t></w:r>' + `
  '<w:r><w:t>1</w:t></w:r>' + `
  '<w:r><w:t>5Kms. The city can be no larger than this size.</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData>' + `
  '</pkg:part>' + `
  '</pkg:package>'

$runRange.InsertXML($rule1Xml)

# ---------------------------------------------------------------------
# Change 2: Add new "Rule 6" text into the (currently empty) placeholder
# paragraph that sits right after the "Rule 5" paragraph.
# ---------------------------------------------------------------------
$rule6Text = "Rule 6: The City should be designed based on scientific " + `
             "principles within social science, empirical evidence " + `
             "should exist for each statement that is made. "

$rule6ParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Trim() -eq "") {
        if ($i -gt 1 -and $d.Paragraphs($i - 1).Range.Text -like "*Rule 5*") {
            $rule6ParaIndex = $i
        }
    }
}

$p6 = $d.Paragraphs($rule6ParaIndex)
$rIns = $d.Range($p6.Range.Start, $p6.Range.Start)
$rIns.Text = $rule6Text

# ---------------------------------------------------------------------
# Change 3: Split the paragraph that holds the closing "For now, ..."
# text together with the "_GoBack" bookmark, so that:
#   - the bookmark stays alone in its own (indented) paragraph
#   - the "For now, ..." text moves into a brand-new paragraph after it
# ---------------------------------------------------------------------
$closingText = "For now, these are the rules that I think are worth " + `
               "giving a thought about, I shall update them as I think " + `
               "becomes necessary. I shall attempt to keep them at a " + `
               "minimum. "

$bookmarkParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*For now,*") {
        $bookmarkParaIndex = $i
    }
}

$pClose = $d.Paragraphs($bookmarkParaIndex)
$closeStart = $pClose.Range.Start
$closeEnd = $closeStart + $closingText.Length

# Remove the visible "For now, ..." text, leaving only the bookmark
# markers behind in the original paragraph.
$delRange = $d.Range($closeStart, $closeEnd)
$delRange.Text = ""

# Insert a brand-new paragraph right after the (now bookmark-only)
# paragraph, containing the closing text.
$pBookmarkOnly = $d.Paragraphs($bookmarkParaIndex)
$insertPoint = $pBookmarkOnly.Range.End - 1
$rAfter = $d.Range($insertPoint, $insertPoint)
$rAfter.InsertAfter("`r" + $closingText)

# Now give the bookmark-only paragraph the same indent formatting used by
# the numbered "Rule" paragraphs (the newly split-off paragraph keeps the
# default / no special indent formatting).
$pBookmarkOnly2 = $d.Paragraphs($bookmarkParaIndex)
$pBookmarkOnly2.Format.LeftIndent = 36
$pBookmarkOnly2.Format.FirstLineIndent = -36
